# PM10 Tidsregistrering for Nikolaj — log two more work entries (row 32 & 33)
# on 2020-03-11 (serial 43901): "SD0701/2/3 og DCD07" (Designer, 09:30-14:00)
# and "UC07 implement af OC07" (implenter, 14:00-15:30). Existing G/H formulas
# (shared formula si="0" and running SUM) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("A32").Value = "SD0701/2/3 og DCD07"
$ws.Range("B32").Value = "Designer"
$ws.Range("C32").Value = 43901
$ws.Range("D32").Value = 0.39583333333333331
$ws.Range("E32").Value = 0.58333333333333337

$ws.Range("A33").Value = "UC07 implement af OC07"
$ws.Range("B33").Value = "implenter"
$ws.Range("C33").Value = 43901
$ws.Range("D33").Value = 0.58333333333333337
$ws.Range("E33").Value = 0.64583333333333337

# Match the author's final selection/scroll position recorded in the sheet view.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A34").Select()
